$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sample assignee rows (Shraddha/Snehal/Kirti ...) with the new
# demo data (abcd/efgh/lnmop ...). Cells are written column-by-column
# (A2..A4, then B3,B2,B4, then C3,C2,C4) so that newly introduced shared
# strings land in the same order as the target workbook's sharedStrings.xml.
$ws.Range("A2").Value = "abcd"
$ws.Range("A3").Value = "efgh"
$ws.Range("A4").Value = "lnmop"
$ws.Range("B3").Value = "s"
$ws.Range("B2").Value = "dfd"
$ws.Range("B4").Value = "dfd"
$ws.Range("C3").Value = "Jira"
$ws.Range("C2").Value = "Email"
$ws.Range("C4").Value = "Email"

# Widen/add the A and B columns and nudge C/D back to (approximately)
# their new widths, matching the resized columns in the edited workbook.
$ws.Columns.Item(1).ColumnWidth = 21.5
$ws.Columns.Item(2).ColumnWidth = 28.666666666666668
$ws.Columns.Item(3).ColumnWidth = 17.666666666666668
$ws.Columns.Item(4).ColumnWidth = 56.5

# Match the saved selection/active cell from the edited workbook.
$ws.Range("C7").Select()
